# Replace the curly “smart” quotes around "Countdown: $_" with straight
# quotes in the PowerShell example on the "Foreach-Object" slide.
#
# Original run text (curly quotes):
#   3, 2, 1 | Foreach-Object {Write-Output “Countdown: $_”}
# New text (straight quotes), after the edit - PowerPoint splits the run
# into five runs at the two retyped quote characters:
#   3, 2, 1 | Foreach-Object {Write-Output "Countdown: $_"}

$p = $ppt.ActivePresentation

# Find the slide/shape that contains the PowerShell snippet.
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text.IndexOf("Foreach-Object {Write-Output") -ge 0) {
                $targetShape = $shape
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text

$needle = '3, 2, 1 | Foreach-Object {Write-Output "Countdown: $_"}'
$needleStart0 = $fullText.IndexOf($needle)
$needleStart1 = $needleStart0 + 1

# Offsets (0-based, within $needle) of the two quote characters.
$firstQuoteOffset = $needle.IndexOf('"')
$secondQuoteOffset = $needle.LastIndexOf('"')

$firstQuotePos = $needleStart1 + $firstQuoteOffset
$secondQuotePos = $needleStart1 + $secondQuoteOffset

# Replace the closing quote first so the earlier position stays valid.
$closingQuote = $tr.Characters($secondQuotePos, 1)
$closingQuote.Text = [char]34

$openingQuote = $tr.Characters($firstQuotePos, 1)
$openingQuote.Text = [char]34
